# The deck currently has its "live" design theme (the one actually wired
# to the slide master / presentation, physically stored as theme2.xml)
# using the "Integral" / "Red Violet" colour scheme, while the unused
# theme part (theme1.xml, referenced only by the notes master) carries
# the default "Office Theme" / "Office" colour scheme.
#
# The target edit swaps the two theme parts' contents, so the design
# that actually renders the deck becomes the plain "Office Theme"
# colours. We reproduce that swap through the exposed COM surface by
# writing the "Office" theme colour values (dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink, in that fixed 1-12 index order) onto the presentation's
# live ThemeColorScheme.

function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

# Office Theme colour scheme, in MsoThemeColorSchemeIndex order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$colorScheme.Colors(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1
$colorScheme.Colors(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1
$colorScheme.Colors(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2
$colorScheme.Colors(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2
$colorScheme.Colors(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1
$colorScheme.Colors(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2
$colorScheme.Colors(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3
$colorScheme.Colors(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4
$colorScheme.Colors(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5
$colorScheme.Colors(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6
$colorScheme.Colors(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink
$colorScheme.Colors(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink
